# Add a UUID field ("line_item_id") to the line_items.tsv section of the data
# dictionary. This inserts one new row right before the existing "metric_id"
# row (row 34), which pushes all subsequent line_items.tsv rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 34; Excel shifts rows 34:51 down to 35:52 and
# carries the existing row formatting/styles down with them.
[void]$ws.Rows(34).Insert()
$ws.Rows(34).RowHeight = 90

# Populate the new row with the line_item_id field definition.
$ws.Range("A34").Value = "line_items.tsv"
$ws.Range("B34").Value = "line_item_id"
$ws.Range("C34").Value = "A universally unique identifier (UUID) for each line item"
$ws.Range("D34").Value = "UUIDs generated using the R id package (https://cran.r-project.org/web/packages/ids/index.html), uuid function"

# Reflect the view/selection state where the author was last working.
[void]$ws.Activate()
[void]$ws.Range("D35").Select()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
